$d = $word.ActiveDocument

# 1. Update the date/time stamp in the title block.
$d.Content.Find.Execute(
    "July   5, 2021 (03:32:31 PM)", $false, $false, $false, $false, $false,
    $true, 1, $false, "July   5, 2021 (03:41:26 PM)", 2) | Out-Null

# 2. Heading text: "Static members in non-static class" -> "...in a non-static class"
$d.Content.Find.Execute(
    "Static members in non-static class", $false, $false, $false, $false, $false,
    $true, 1, $false, "Static members in a non-static class", 2) | Out-Null

# 3. Body text: "both static or non-static" -> "both static and non-static"
$d.Content.Find.Execute(
    "A non-static class can contain both static or non-static class members.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "A non-static class can contain both static and non-static class members.", 2) | Out-Null

# 4. Comment text inside code block.
$d.Content.Find.Execute(
    "// does this work? uncomment next line to check!", $false, $false, $false, $false, $false,
    $true, 1, $false, "// does this work? uncomment to check", 2) | Out-Null

Write-Host "done"
